$d = $word.ActiveDocument

# Locate the "Requisitos" entry paragraph that must stay, and the
# trailing "(c) 2020 ..." footer paragraph that must go. Everything from
# right after the former up to (and including) the end of the latter is
# removed: the blank spacer paragraph, the "Ver no Jupiter..." paragraph,
# and the "(c) 2020 ..." paragraph itself.
$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text
    if ($text -like "*LOM3015: Termodinâmica de Materiais (Requisito fraco)*") {
        $startPara = $p
    }
    if ($text -like "*Contact: luizeleno@usp.br*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $deleteRange = $d.Range($startPara.Range.End, $endPara.Range.End)
    $deleteRange.Delete()
}
